$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.339.49"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "1.867.43"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3964"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08023"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9983"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "1.876.70"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.009"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.237"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001042"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06634"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "28.361.40"
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.466"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.266"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "2.094.80"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.123"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.488"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9657"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09485"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.350"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.367"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06082"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02246"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.352"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.186"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5940"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1872"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.290"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5575"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.954"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06856"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.036"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.41%  "
